$d = $word.ActiveDocument

# Locate the complex field (the M2Doc template expression currently stored
# as a real Word field: fldChar begin / instrText* / fldChar end) and the
# paragraph that contains it.
$f = $null
for ($i = 1; $i -le $d.Fields.Count; $i++) {
    $candidate = $d.Fields.Item($i)
    if ($candidate.Code.Text -match "asImageByRepresentationDescriptionName") {
        $f = $candidate
        break
    }
}
if ($f -eq $null) {
    throw "Could not find the asImageByRepresentationDescriptionName field"
}

$paraIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($f.Code.Start -ge $p.Range.Start -and $f.Code.Start -lt $p.Range.End) {
        $paraIndex = $i
        break
    }
}
if ($paraIndex -eq -1) {
    throw "Could not find the paragraph containing the field"
}

$para = $d.Paragraphs.Item($paraIndex)
$full = $para.Range
# Trim off the trailing paragraph mark so the paragraph's own <w:p>/<w:pPr>
# survive untouched and only its run content gets rewritten.
$target = $d.Range($full.Start, $full.End - 1)

# Rebuild the field's runs as plain literal text runs: the begin/end
# fldChar runs are dropped, each instrText run becomes a <w:t> run with the
# same text, and the whole expression is wrapped in "{" ... "}" (M2Doc's
# own field-rewriter syntax) on the first/last run respectively.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidP="00340939" w:rsidR="00A10D75" w:rsidRDefault="00474E78"><w:pPr><w:widowControl w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>{m:</w:t></w:r><w:r w:rsidR="003141BA"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>self</w:t></w:r><w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="004C24F3" w:rsidRPr="004C24F3"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>asImageByRepresentationDescriptionName</w:t></w:r><w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>(</w:t></w:r><w:r w:rsidR="00327A56"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&apos;</w:t></w:r><w:r w:rsidR="00AE2CDB"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Entities</w:t></w:r><w:r w:rsidR="003141BA"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>&apos;</w:t></w:r><w:r w:rsidR="00C53443"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r><w:r w:rsidR="0014442D"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>-&gt;first()</w:t></w:r><w:r w:rsidR="000E5422"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.setWidth(300)}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml)
